$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 7 rows before row 20 to make room for the extra account-statement rows ---
$ws.Rows("20:26").Insert()

# Row 26 (the new last data row) takes over the special "final row" style currently on row 19
$ws.Range("B19:J19").Copy()
$ws.Range("B26:J26").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 19 reverts to the regular body-row style (matching row 18); rows 20-25 get it too
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Header values ---
$ws.Range("E11").Value = 587990
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 8

# --- Data table rows 16-26 ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73124030"
$ws.Range("D16").Value = "RICHARD BRIEVA QUINTANA"
$ws.Range("E16").Value = "2502"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 0

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "19774068"
$ws.Range("D17").Value = "TONIS ALBERTO BASTIDAS CARVAJAL"
$ws.Range("E17").Value = "2501"
$ws.Range("F17").Value = 52000
$ws.Range("G17").Value = 1300000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "19774068"
$ws.Range("D18").Value = "TONIS ALBERTO BASTIDAS CARVAJAL"
$ws.Range("E18").Value = "2412"
$ws.Range("F18").Value = 52000
$ws.Range("G18").Value = 1300000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "92400506"
$ws.Range("D19").Value = "ALBERTO JOSE BLANCO JULIO"
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "92400506"
$ws.Range("D20").Value = "ALBERTO JOSE BLANCO JULIO"
$ws.Range("E20").Value = "2506"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "92400506"
$ws.Range("D21").Value = "ALBERTO JOSE BLANCO JULIO"
$ws.Range("E21").Value = "2505"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "92400506"
$ws.Range("D22").Value = "ALBERTO JOSE BLANCO JULIO"
$ws.Range("E22").Value = "2504"
$ws.Range("F22").Value = 56940
$ws.Range("G22").Value = 1423500

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "92400506"
$ws.Range("D23").Value = "ALBERTO JOSE BLANCO JULIO"
$ws.Range("E23").Value = "2503"
$ws.Range("F23").Value = 56940
$ws.Range("G23").Value = 1423500

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "92400506"
$ws.Range("D24").Value = "ALBERTO JOSE BLANCO JULIO"
$ws.Range("E24").Value = "2502"
$ws.Range("F24").Value = 47450
$ws.Range("G24").Value = 1423500

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1047393044"
$ws.Range("D25").Value = "SABEL TILVEZ FERNANDEZ"
$ws.Range("E25").Value = "2502"
$ws.Range("F25").Value = 47450
$ws.Range("G25").Value = 877803

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1043649211"
$ws.Range("D26").Value = "VICTOR MANUEL SIMANCAS SARMIENTO"
$ws.Range("E26").Value = "2502"
$ws.Range("F26").Value = 47450
$ws.Range("G26").Value = 1423500

# --- Column D width widened to fit the longer names ---
$ws.Columns("D").ColumnWidth = 36

Write-Host "done"
